# The source workbook stores B2 as a shared-string "75939" (numeric-looking
# text) rather than a real number. Excel's normal type-inference would turn a
# bare numeric literal typed into a General-formatted cell into a Number, so
# force Text formatting just long enough to commit the new literal as a
# string, then restore the cell's style back to Normal/General so the net
# effect is a pure value change (same type "s", same style) matching the
# target edit: B2 75939 -> 1035.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B2")
$cell.NumberFormat = "@"
$cell.Value = "1035"
$cell.Style = "Normal"
